# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de handback
# has completed (mirrors the zh-cn handback that was already recorded), and
# that the Overview status text for both languages is now "Handed back: in
# sync with en-US" instead of the stale "Ready for handoff".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: the "zh-cn" / "de-de" status cells move from
# "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Columns E/F grow to fit the new, longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet: fill in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns (I/J/K) for both rows, now that the
# handback has been recorded.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("J2").Value = "0994a4d9-84b7-4981-a007-5742534ba6d3.ff2e8eb7d28de7923dc04add43e747c214b3a421.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-12 18:34:37"

$zhcn.Range("J3").Value = "907d2966-b6bb-49eb-81a2-42d2029526fb.b8e254f6b013de04d958da23b32004e52810808f.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-12 18:34:37"

# Rebuild the hyperlinks collection so the pre-existing "Latest Handoff File"
# links (column A) are joined by new "Latest Target File" links (column I)
# that point at the same source document.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3cd8f90253581bcb0cd2b3c8233ea36600b44a07/e2e/0994a4d9-84b7-4981-a007-5742534ba6d3.md", "", "", "0994a4d9-84b7-4981-a007-5742534ba6d3.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3cd8f90253581bcb0cd2b3c8233ea36600b44a07/e2e/0994a4d9-84b7-4981-a007-5742534ba6d3.md", "", "", "0994a4d9-84b7-4981-a007-5742534ba6d3.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3cd8f90253581bcb0cd2b3c8233ea36600b44a07/e2e/907d2966-b6bb-49eb-81a2-42d2029526fb.md", "", "", "907d2966-b6bb-49eb-81a2-42d2029526fb.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3cd8f90253581bcb0cd2b3c8233ea36600b44a07/e2e/907d2966-b6bb-49eb-81a2-42d2029526fb.md", "", "", "907d2966-b6bb-49eb-81a2-42d2029526fb.md")

# Columns C (Status), I (Latest Target File) and J (Latest Handback File)
# widen to fit the longer file-name/status content now shown.
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet: same treatment as zh-cn, but the handback datetime is its own
# distinct value (2016-08-12 18:34:46) rather than the zh-cn one.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("J2").Value = "0994a4d9-84b7-4981-a007-5742534ba6d3.ff2e8eb7d28de7923dc04add43e747c214b3a421.de-de.xlf"
$dede.Range("K2").Value = "2016-08-12 18:34:46"

$dede.Range("J3").Value = "907d2966-b6bb-49eb-81a2-42d2029526fb.b8e254f6b013de04d958da23b32004e52810808f.de-de.xlf"
$dede.Range("K3").Value = "2016-08-12 18:34:46"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3cd8f90253581bcb0cd2b3c8233ea36600b44a07/e2e/0994a4d9-84b7-4981-a007-5742534ba6d3.md", "", "", "0994a4d9-84b7-4981-a007-5742534ba6d3.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3cd8f90253581bcb0cd2b3c8233ea36600b44a07/e2e/0994a4d9-84b7-4981-a007-5742534ba6d3.md", "", "", "0994a4d9-84b7-4981-a007-5742534ba6d3.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3cd8f90253581bcb0cd2b3c8233ea36600b44a07/e2e/907d2966-b6bb-49eb-81a2-42d2029526fb.md", "", "", "907d2966-b6bb-49eb-81a2-42d2029526fb.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/3cd8f90253581bcb0cd2b3c8233ea36600b44a07/e2e/907d2966-b6bb-49eb-81a2-42d2029526fb.md", "", "", "907d2966-b6bb-49eb-81a2-42d2029526fb.md")

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated"
